$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the price/volume cells to remain plain text (matches the source
# data, which stores these as text strings like "302.42" / "-0.58%"),
# instead of Excel auto-converting numeric-looking text into numbers.
$textCells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E21","D22","E22","D23","E23","E24","D25","E25","D26","E26","D39","E39","D40","E40","E41","D42","E42","D43","E43","D44","E44","D45","E45","D46","E46","D47","E47","D48","E48","D49","E49","D50","E50","D51","E51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Row re-ranking (rows 8-17): Coin name + Link columns ---
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"

# --- Updated Price (D) / Volume 1h (E) values ---
$ws.Range("D2").Value = "302.42"
$ws.Range("E2").Value = "-0.58%"
$ws.Range("D3").Value = "36.68"
$ws.Range("E3").Value = "3.17%"
$ws.Range("D4").Value = "4.982"
$ws.Range("E4").Value = "-1.92%"
$ws.Range("D5").Value = "0.07706"
$ws.Range("E5").Value = "-1.47%"
$ws.Range("D6").Value = "2.106"
$ws.Range("E6").Value = "-6.96%"
$ws.Range("D7").Value = "7.946"
$ws.Range("E7").Value = "-1.68%"
$ws.Range("D8").Value = "4.020"
$ws.Range("E8").Value = "-0.69%"
$ws.Range("D9").Value = "0.9128"
$ws.Range("E9").Value = "-1.71%"
$ws.Range("D10").Value = "0.09667"
$ws.Range("E10").Value = "2.97%"
$ws.Range("D11").Value = "0.1848"
$ws.Range("E11").Value = "0.75%"
$ws.Range("D12").Value = "0.08549"
$ws.Range("E12").Value = "-0.11%"
$ws.Range("D13").Value = "0.03553"
$ws.Range("E13").Value = "-1.89%"
$ws.Range("D14").Value = "0.09954"
$ws.Range("E14").Value = "-0.18%"
$ws.Range("D15").Value = "0.001466"
$ws.Range("E15").Value = "-1.03%"
$ws.Range("D16").Value = "0.005678"
$ws.Range("E16").Value = "0.57%"
$ws.Range("D17").Value = "3.465"
$ws.Range("E17").Value = "-0.34%"
$ws.Range("D18").Value = "2.219"
$ws.Range("E18").Value = "7.32%"
$ws.Range("D19").Value = "0.3387"
$ws.Range("E19").Value = "-0.57%"
$ws.Range("D20").Value = "0.1329"
$ws.Range("E20").Value = "0.52%"
$ws.Range("D21").Value = "4.754"
$ws.Range("E21").Value = "4.03%"
$ws.Range("D22").Value = "0.2196"
$ws.Range("E22").Value = "-1.81%"
$ws.Range("D23").Value = "0.04581"
$ws.Range("E23").Value = "-2.06%"
$ws.Range("E24").Value = "12.22%"
$ws.Range("D25").Value = "0.001228"
$ws.Range("E25").Value = "-0.34%"
$ws.Range("D26").Value = "0.0001397"
$ws.Range("E26").Value = "7.40%"
$ws.Range("D39").Value = "0.01762"
$ws.Range("E39").Value = "-0.98%"
$ws.Range("D40").Value = "0.04617"
$ws.Range("E40").Value = "-2.22%"
$ws.Range("E41").Value = "-6.11%"
$ws.Range("D42").Value = "0.1390"
$ws.Range("E42").Value = "-2.19%"
$ws.Range("D43").Value = "0.007640"
$ws.Range("E43").Value = "-4.46%"
$ws.Range("D44").Value = "0.002156"
$ws.Range("E44").Value = "-5.91%"
$ws.Range("D45").Value = "0.01035"
$ws.Range("E45").Value = "14.10%"
$ws.Range("D46").Value = "0.00006298"
$ws.Range("E46").Value = "1.67%"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").Value = "-0.25%"
$ws.Range("D48").Value = "0.0005809"
$ws.Range("E48").Value = "0.15%"
$ws.Range("D49").Value = "37.16"
$ws.Range("E49").Value = "594.10%"
$ws.Range("D50").Value = "0.001997"
$ws.Range("E50").Value = "-25.83%"
$ws.Range("D51").Value = "0.00002096"
$ws.Range("E51").Value = "-0.25%"
